$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B (StatQuery query/result column),
# shifting the old B (Neo4jData) and C (WebData query) columns right to C and D.
$ws.Range("B1").EntireColumn.Insert() | Out-Null

# Match new column B's width to column A's width (same as source column copied).
$ws.Columns.Item(2).ColumnWidth = $ws.Columns.Item(1).ColumnWidth

# Header for the new column.
$ws.Range("B1").Value = "StatQuery"

# New stat-bar query content, word-wrapped like A2.
$ws.Range("B2").Value = 'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN[''Rottweiler'']  OPTIONAL MATCH (f:file)-[*]->(c), (samp:sample)-[*]->(c) WITH DISTINCT c AS c, p, s, demo, diag, f, samp RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(samp)) as number_of_sample , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(s.clinical_study_designation)) as number_of_study'
$ws.Range("B2").WrapText = $true

# Update the active selection to A2 (matches the saved view state).
$ws.Range("A2").Select() | Out-Null
